$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count"); remaining columns F:K shift left to E:J
$ws.Range("E:E").Delete()
